# Fruta / hortaliza, semanal
# Insert three new weekly price rows (Royal Glory, Region de O'Higgins / Paine)
# into the "Durazno" price table, shifting the rows below each insertion point
# down, matching the new weekly data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that are constant for every data row in this sheet
$colA = 10
$colB = "Vega Modelo de Temuco"
$colC = "La Araucanía"
$colE = 9
$colF = "Fruta"
$colG = 100103
$colH = "Frutos de hueso (carozo)"
$colI = 100103004
$colJ = "Durazno"

function Set-DataRow($RowNum, $D, $K, $L, $M, $N, $O, $P, $Q, $R, $S, $T) {
    $ws.Cells.Item($RowNum, 1).Value = $colA
    $ws.Cells.Item($RowNum, 2).Value = $colB
    $ws.Cells.Item($RowNum, 3).Value = $colC
    $ws.Cells.Item($RowNum, 4).Value = $D
    $ws.Cells.Item($RowNum, 5).Value = $colE
    $ws.Cells.Item($RowNum, 6).Value = $colF
    $ws.Cells.Item($RowNum, 7).Value = $colG
    $ws.Cells.Item($RowNum, 8).Value = $colH
    $ws.Cells.Item($RowNum, 9).Value = $colI
    $ws.Cells.Item($RowNum, 10).Value = $colJ
    $ws.Cells.Item($RowNum, 11).Value = $K
    $ws.Cells.Item($RowNum, 12).Value = $L
    $ws.Cells.Item($RowNum, 13).Value = $M
    $ws.Cells.Item($RowNum, 14).Value = $N
    $ws.Cells.Item($RowNum, 15).Value = $O
    $ws.Cells.Item($RowNum, 16).Value = $P
    $ws.Cells.Item($RowNum, 17).Value = $Q
    $ws.Cells.Item($RowNum, 18).Value = $R
    $ws.Cells.Item($RowNum, 19).Value = $S
    $ws.Cells.Item($RowNum, 20).Value = $T
}

# --- Insert two new rows at 130:131 (Royal Glory, 2022-01-06) ---
$ws.Rows("130:131").Insert()

Set-DataRow 130 44567 "Royal Glory" "Primera" 250 16000 16000 16000 "`$/bandeja 18 kilos granel" "Región de O'Higgins" 889 18
Set-DataRow 131 44567 "Royal Glory" "Primera" 8 320000 320000 320000 "`$/bins (400 kilos)" "Región de O'Higgins" 800 400

# --- Insert one new row at 202 (Royal Glory, 2022-01-07, Paine) ---
$ws.Rows("202:202").Insert()

Set-DataRow 202 44568 "Royal Glory" "Primera" 125 17000 17000 17000 "`$/bandeja 18 kilos granel" "Paine" 944 18
